# Cleaning Samples Template - sample list checker update
#
# - Clears the "Other Peak(s)" value in G5 (was a bogus "string to fail
#   here" placeholder) and the stray "sample  " value in B7, so the
#   "other peaks" parser only sees real (rt, conc) tuple strings.
# - Formats the Sample ID column (A3:A12) as Text ("@") so sample IDs
#   keep any leading zeros / aren't coerced to numbers.
# - Moves the active selection to G5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the bad placeholder values.
$ws.Range("G5").Value = ""
$ws.Range("B7").Value = ""

# Sample ID column -> Text format for all data rows.
$ws.Range("A3:A12").NumberFormat = "@"

# Update the selected cell/active cell.
[void]$ws.Range("G5").Select()
